# Actualización desde MV -datos-
# Append the new daily rows (21-09-2021 .. 01-10-2021) to the bottom of the
# "14 días 2021 - Diaria" table, mirroring the existing row layout
# (columns A..G: Serie, Cupo, Monto demandado, Monto total adjudicado,
# Monto adjudicado bancos, Monto adjudicado AFP, Tasa).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("21-09-2021", 2600000, 2703000, 2600000, 1337000, 1263000, 1.48),
    @("22-09-2021", 2600000, 4473000, 3900000, 2607000, 1293000, 1.48),
    @("23-09-2021", 1000000, 2998000, 1500000,  867000,  633000, 1.45),
    @("24-09-2021",  700000, 2168000, 1050000,  277000,  773000, 1.42),
    @("27-09-2021", 1000000, 2633000, 1500000, 1025000,  475000, 1.4),
    @("28-09-2021", 1000000, 2358000, 1500000,  937000,  563000, 1.45),
    @("29-09-2021", 1500000, 2335000, 2250000, 1975000,  275000, 1.47),
    @("30-09-2021", 3000000, 3433000, 3000000, 1967000, 1033000, 1.47),
    @("01-10-2021", 4000000, 5270000, 4000000, 3213000,  787000, 1.47)
)

$startRow = 182
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]

    # Column A holds the date label as text. Excel's automatic data-type
    # detection can reinterpret an unambiguous-looking dd-mm-yyyy string
    # (e.g. "01-10-2021", where both parts are <= 12) as a date serial, so
    # force the cell to Text before writing, then drop the explicit format
    # again so the cell ends up unstyled, matching the rest of the column.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0]
    $cellA.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
}
